# Applies the "Отчет6.docx" edits:
#  1. Executor name: Шаронов Никита Андреевич -> Медведев Владислав Александрович
#     (split into ": " + the new name, two runs)
#  2. Group number: 207 -> 205 (split into ": 20" + "5", two runs)
#  3. Merge the two runs that make up " _______________________" into one run
#  4. Mark the five inline drawings as NoProofing (-> <w:noProof/> in rPr)
#  5. Insert <w:lastRenderedPageBreak/> before "Задание 2:".."Задание 5:"
#  6. Remove the (now stale) _GoBack bookmark
#
# NOTE: this engine's Range.InsertXML only behaves predictably when the
# target Range spans a *whole paragraph's text* (start of first run's text
# to end of last run's text) - sub-run-range InsertXML calls corrupt the
# surrounding paragraph. So every run-split edit below first Finds the
# whole paragraph text, then resubmits the complete <w:p> as replacement
# XML (the engine keeps the original <w:p>'s own rsid/attributes).

$d   = $word.ActiveDocument
$sel = $word.Selection
$wNs = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

function Replace-ParagraphXml([string]$findText, [string]$pInnerXml) {
    $sel.Find.Execute($findText, $true, $false, $false, $false, $false, $true, 1, $false, $null, 0) | Out-Null
    $r = $d.Range($sel.Start, $sel.End)
    $r.InsertXML('<w:p ' + $wNs + '>' + $pInnerXml + '</w:p>')
}

# 1. Executor name ------------------------------------------------------
Replace-ParagraphXml "Исполнитель: Шаронов Никита Андреевич" (
    '<w:pPr><w:spacing w:line="360" w:lineRule="auto"/><w:rPr><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr></w:pPr>' +
    '<w:r><w:rPr><w:b/><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr><w:t>Исполнитель</w:t></w:r>' +
    '<w:r><w:rPr><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr><w:t xml:space="preserve">: </w:t></w:r>' +
    '<w:r><w:rPr><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr><w:t>Медведев Владислав Александрович</w:t></w:r>'
)

# 2. Group number ---------------------------------------------------------
Replace-ParagraphXml "Группа: 207" (
    '<w:pPr><w:spacing w:line="360" w:lineRule="auto"/><w:rPr><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr></w:pPr>' +
    '<w:r><w:rPr><w:b/><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr><w:t>Группа</w:t></w:r>' +
    '<w:r><w:rPr><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr><w:t>: 20</w:t></w:r>' +
    '<w:r><w:rPr><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr><w:t>5</w:t></w:r>'
)

# 3. Merge the signature-line underscore runs ------------------------------
Replace-ParagraphXml "Подпись преподавателя   _______________________" (
    '<w:pPr><w:spacing w:line="360" w:lineRule="auto"/><w:rPr><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr></w:pPr>' +
    '<w:r><w:rPr><w:b/><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr><w:t xml:space="preserve">Подпись преподавателя  </w:t></w:r>' +
    '<w:r><w:rPr><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr><w:t xml:space="preserve"> _______________________</w:t></w:r>'
)

# 4 & 5. Walk the five task drawings --------------------------------------
$taskNames = @("Задание 1:", "Задание 2:", "Задание 3:", "Задание 4:", "Задание 5:")
for ($i = 1; $i -le 5; $i++) {
    $shp = $d.InlineShapes.Item($i)
    $shp.Range.NoProofing = $true

    if ($i -ge 2) {
        $label = $taskNames[$i - 1]
        Replace-ParagraphXml $label (
            '<w:pPr><w:ind w:firstLine="709"/><w:jc w:val="both"/><w:rPr><w:b/><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr></w:pPr>' +
            '<w:r><w:rPr><w:b/><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr><w:lastRenderedPageBreak/><w:t>' + $label + '</w:t></w:r>'
        )
    }
}

# 6. Drop the stale _GoBack bookmark ---------------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}
